$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder "Sierra Leona" to appear right after "Aruba" (before "Togo"),
# pushing "Togo" and "Monaco" down by one row, and refresh the daily
# case numbers for these three countries (rows 149-151).
$ws.Range("A149").Value = "Sierra Leona"
$ws.Range("B149").Value = 99
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 10
$ws.Range("E149").Value = 85
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 4

$ws.Range("A150").Value = "Togo"
$ws.Range("B150").Value = 99
$ws.Range("C150").Value = 1
$ws.Range("D150").Value = 62
$ws.Range("E150").Value = 31
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 6

$ws.Range("A151").Value = "Monaco"
$ws.Range("B151").Value = 95
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 42
$ws.Range("E151").Value = 49
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 4

# Update the "last updated" timestamp string in A1.
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 05:52"
